# Reverses the order of the two comma-separated entries in the "Recorded By"
# column (G) for the specific rows affected by the sync, e.g.
#   "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
#   "backup@backdoor.com, System" -> "System, backup@backdoor.com"
#   "dnasr281@gmail.com, admin@admin.com" -> "admin@admin.com, dnasr281@gmail.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$rows = @(4,10,11,12,13,14,15,17,18,19,20,21,22,24,26,30,36,37,38,39,40,41,43,44,45,46,47,48,50,52,56,62,63,64,65,66,67,69,70,71,72,73,74,76,78,83,84,85,86,87,90,92,93,94,96,99,101,109,110,111,112,113,116,118,119,120,122,125,127,135,136,137,138,139,142,144,145,146,148,151,153)

foreach ($r in $rows) {
    $cell = $ws.Range("G$r")
    $current = [string]$cell.Value2
    $parts = $current.Split(",")
    if ($parts.Length -eq 2) {
        $first = $parts[0].Trim()
        $second = $parts[1].Trim()
        $cell.Value2 = "$second, $first"
    }
}

$wb.Save()
